$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4085.191239257138
$ws.Range("C3").Value = 4085.191239257138
$ws.Range("C4").Value = 3882.038292339431
$ws.Range("C5").Value = 3882.038292339431
$ws.Range("C6").Value = 3751.331342649916
$ws.Range("C7").Value = 3741.397770715599
$ws.Range("C8").Value = 3741.397770715599
$ws.Range("C9").Value = 3741.397770715599
$ws.Range("C10").Value = 3582.818119354935
$ws.Range("C11").Value = 3582.818119354935
$ws.Range("C12").Value = 3582.818119354935
